$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "test201"
$ws.Range("B2").Value = 23071223
$ws.Range("C2").Value = "narendra658"
$ws.Range("D2").Value = "S#w3Kk7%"
